# Extend the macro correlation table from 3 columns (mean:var, mean:iqr,
# mean:rvar) to 7 columns, adding mean:skew and the median:* variants, and
# refresh all correlation values (separating labor-market vs stock-market
# macro analysis; only the needed results are exported).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (B1:H1) -------------------------------------------------
$headers = @("mean:var", "mean:iqr", "mean:rvar", "mean:skew", "median:var", "median:iqr", "median:rvar")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Formula = '="' + $headers[$i] + '"'
}

# ---- Data rows (row 2..8, columns B..H) ---------------------------------
$data = @(
    @("-0.28**",  "-0.36***", "-0.43***", "0.11",   "-0.32***", "-0.29**", "-0.49***"),
    @("-0.32***", "-0.38***", "-0.29***", "0.38***", "-0.28**", "-0.27**", "-0.35***"),
    @("-0.31***", "-0.32***", "-0.31***", "0.25**",  "-0.31***", "-0.26**", "-0.41***"),
    @("-0.26**",  "-0.28**",  "-0.26**",  "0.13",    "-0.24**",  "-0.21*",  "-0.32***"),
    @("-0.19",    "-0.26**",  "-0.21*",   "0.05",    "-0.33***", "-0.31***", "-0.41***"),
    @("-0.33***", "-0.41***", "-0.4***",  "0.15",    "-0.4***",  "-0.38***", "-0.48***"),
    @("-0.46***", "-0.51***", "-0.37***", "0.22*",   "-0.43***", "-0.43***", "-0.47***")
)

$dataCols = @("B", "C", "D", "E", "F", "G", "H")
for ($r = 0; $r -lt $data.Length; $r++) {
    $excelRow = $r + 2
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $dataCols.Length; $c++) {
        $ws.Range($dataCols[$c] + $excelRow).Formula = '="' + $rowVals[$c] + '"'
    }
}

# Convert all the helper formulas above into literal text values, so the
# cells are plain strings (like the rest of the sheet) rather than live
# formulas. Doing the whole block in one Copy/PasteSpecial keeps the
# workbook's style table untouched.
$full = $ws.Range("B1:H8")
$full.Copy()
$full.PasteSpecial(-4163)  # xlPasteValues

# ---- Formatting: extend the bold/bordered/centered header style and the
# bold column-A style to the new columns/rows (E1:H1 header style; keep
# existing column A row styles as-is since A only spans rows 2-8 already).
$ws.Range("B1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
